$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168; this shifts the existing rows 168..223
# down to 169..224 (values, not just formatting), matching the target diff.
$ws.Rows(168).Insert()

# Populate the newly inserted row 168 with the new record.
$ws.Cells.Item(168, 1).Value = 11
$ws.Cells.Item(168, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(168, 3).Value = "Bíobío"
$ws.Cells.Item(168, 4).Value = 45135
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 6).Value = 100112021
$ws.Cells.Item(168, 7).Value = "Ají"
$ws.Cells.Item(168, 8).Value = "Americana (o)"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 40
$ws.Cells.Item(168, 11).Value = 35000
$ws.Cells.Item(168, 12).Value = 36000
$ws.Cells.Item(168, 13).Value = 35500
$ws.Cells.Item(168, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(168, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(168, 16).Value = 1420
$ws.Cells.Item(168, 17).Value = 25
$ws.Cells.Item(168, 18).Value = "Hortaliza"
